$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A column labels (rows 3-26); A2 unchanged ("model_2_3_0")
$ws.Range("A3").Value = "model_2_3_22"
$ws.Range("A4").Value = "model_2_3_21"
$ws.Range("A5").Value = "model_2_3_20"
$ws.Range("A6").Value = "model_2_3_19"
$ws.Range("A7").Value = "model_2_3_18"
$ws.Range("A8").Value = "model_2_3_17"
$ws.Range("A9").Value = "model_2_3_16"
$ws.Range("A10").Value = "model_2_3_15"
$ws.Range("A11").Value = "model_2_3_14"
$ws.Range("A12").Value = "model_2_3_13"
$ws.Range("A13").Value = "model_2_3_23"
$ws.Range("A14").Value = "model_2_3_12"
$ws.Range("A15").Value = "model_2_3_10"
$ws.Range("A16").Value = "model_2_3_9"
$ws.Range("A17").Value = "model_2_3_8"
$ws.Range("A18").Value = "model_2_3_7"
$ws.Range("A19").Value = "model_2_3_6"
$ws.Range("A20").Value = "model_2_3_5"
$ws.Range("A21").Value = "model_2_3_4"
$ws.Range("A22").Value = "model_2_3_3"
$ws.Range("A23").Value = "model_2_3_2"
$ws.Range("A24").Value = "model_2_3_1"
$ws.Range("A25").Value = "model_2_3_11"
$ws.Range("A26").Value = "model_2_3_24"

# Update B:I columns (rows 2-26) to the new constant metric values
$colVal_B = -0.001236788884735551
$colVal_C = 0.2833859744320865
$colVal_D = -0.1405705737399656
$colVal_E = 0.01893983943866462
$colVal_F = 1.108074069023132
$colVal_G = 0.273603230714798
$colVal_H = 1.107705950737
$colVal_I = 0.6661221385002136

for ($r = 2; $r -le 26; $r++) {
    $ws.Range("B" + $r).Value = $colVal_B
    $ws.Range("C" + $r).Value = $colVal_C
    $ws.Range("D" + $r).Value = $colVal_D
    $ws.Range("E" + $r).Value = $colVal_E
    $ws.Range("F" + $r).Value = $colVal_F
    $ws.Range("G" + $r).Value = $colVal_G
    $ws.Range("H" + $r).Value = $colVal_H
    $ws.Range("I" + $r).Value = $colVal_I
}
